# logboek Rens.xlsx - "temp en webserver changes"
# Adds a new logboek entry (row 9) describing that the gateway can request
# data from the arduino, normalizes the time formatting of the three most
# recent entries (HH:MM:SS -> HH:MM) and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-8 used a HH:MM:SS number format (style index 5); align them with
# the HH:MM format already used elsewhere in the sheet (style index 3).
$ws.Range("D6:E8").NumberFormat = "HH:MM"

# New log entry on row 9.
$ws.Cells.Item(9, 1).Value = "alleen"
$ws.Cells.Item(9, 1).NumberFormat = "General"

$ws.Cells.Item(9, 2).Value = "gateway kan de data van de arduino opvragen"
$ws.Cells.Item(9, 2).NumberFormat = "General"

$ws.Cells.Item(9, 3).Value = 43398
$ws.Cells.Item(9, 3).NumberFormat = "DD/MM/YY"

$ws.Cells.Item(9, 4).Value = 0.375
$ws.Cells.Item(9, 4).NumberFormat = "HH:MM"

$ws.Cells.Item(9, 5).Value = 0.472222222222222
$ws.Cells.Item(9, 5).NumberFormat = "HH:MM"

# Move the active selection to H20, as in the final workbook.
$ws.Range("H20").Select()
